# Adds a new "URL" column (F) to Hoja1 with a quick-reference hyperlink in F2,
# per commit: "se agrega una columna en el archivo de excel con el link para
# consulta rapida" (a column with a link for quick lookup is added).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New header + hyperlink display text (shared strings: "URL", "ver aqui")
$ws.Range("F1").Value = "URL"
$ws.Range("F2").Value = "ver aqui"

# Turn F2 into a real hyperlink (creates the Hyperlink style/font + the
# external relationship automatically).
$ws.Hyperlinks.Add($ws.Range("F2"), "https://example.com/consulta-rapida") | Out-Null

# Column widths to fit the new contents (matches bestFit-style sizing).
$ws.Range("E1").ColumnWidth = 82.375
$ws.Range("F1").ColumnWidth = 7.125

# View tweaks: zoom out a bit and move the active selection.
$ws.Activate()
$excel.ActiveWindow.Zoom = 80
$ws.Range("B3").Select() | Out-Null
